$d = $word.ActiveDocument

# Locate a single "E" character elsewhere in the document that already has the
# exact run formatting we need (rFonts GeoSlab703 MdCn BT, bold, smallCaps,
# themeColor text1 / themeTint d9 / val 262626, single underline) -- the
# "8.- OBSERVACION GENERAL" heading uses identical formatting to our target
# heading and conveniently contains the letter "E" we need.
$srcRng = $d.Content
$srcRng.Find.Execute("OBSERVACION GENERAL", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$eSrc = $d.Range($srcRng.Start + 13, $srcRng.Start + 14)

# Locate the typo "GENRAL" inside the "9.- RECOMENDACION GENRAL" heading.
$target = $d.Content
$target.Find.Execute("GENRAL", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$genStart = $target.Start

# Insert a correctly-formatted "E" between "GEN" and "RAL" (at offset 3 of
# "GENRAL") by copying the formatted run from the source above. Using
# FormattedText (a structural run copy) keeps it as its own run instead of
# being silently re-merged into its neighbours.
$insertionPoint = $d.Range($genStart + 3, $genStart + 3)
$insertionPoint.FormattedText = $eSrc
